# Update countries & provincias Spain
# - Re-sort a handful of countries whose case counts changed position
#   (their shared-string label moves to a different row, and the
#   surrounding numeric columns shift accordingly).
# - Bump the "Datos actualizados" timestamp from 06:20 to 06:50.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp (row 1) ---------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 30 de Marzo de 2020 a las 06:50"

# --- Country name relabels (values only move between existing rows) -----
$ws.Range("A46").Value = "Mexico"
$ws.Range("A47").Value = "Panama"
$ws.Range("A48").Value = "Republica Dominicana"
$ws.Range("A49").Value = "Peru"

$ws.Range("A104").Value = "Sri Lanka"
$ws.Range("A105").Value = "Venezuela"

$ws.Range("A140").Value = "El Salvador"
$ws.Range("A141").Value = "Togo"
$ws.Range("A142").Value = "Zambia"

$ws.Range("A154").Value = "Bahamas"
$ws.Range("A155").Value = "Tanzania"

$ws.Range("A160").Value = "San Martin (Parte Francesa)"
$ws.Range("A161").Value = "Dominica"

# --- Updated statistics for the rows whose underlying country changed ---
function Set-RowStats($row, $b, $c, $d, $e, $f, $g, $h) {
    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 6).Value = $f
    $ws.Cells.Item($row, 7).Value = $g
    $ws.Cells.Item($row, 8).Value = $h
}

Set-RowStats 46 993 145 35 938 1 4 20
Set-RowStats 47 989 0 4 961 36 0 24
Set-RowStats 48 859 0 3 817 0 0 39
Set-RowStats 49 852 0 16 818 40 0 18

Set-RowStats 104 120 3 11 108 5 0 1
Set-RowStats 105 119 0 39 77 6 0 3

Set-RowStats 140 30 6 0 30 0 0 0
Set-RowStats 141 30 5 1 28 0 0 1
Set-RowStats 142 29 0 0 29 0 0 0
